# Apply the edits described by the diff:
# - Update student IDs in A6:A12 from 2025210xxx to 2025310xxx
#   (B6:B12 contain shared formulas referencing column A, so they update automatically)
# - Update the active selection on Sheet1 to C10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = 2025310025
$ws.Range("A7").Value = 2025310026
$ws.Range("A8").Value = 2025310029
$ws.Range("A9").Value = 2025310030
$ws.Range("A10").Value = 2025310032
$ws.Range("A11").Value = 2025310033
$ws.Range("A12").Value = 2025310034

$ws.Activate() | Out-Null
$ws.Range("C10").Select() | Out-Null
